# This workbook records weekly fruit/vegetable price observations.
# A new weekly observation row needs to be inserted right before the
# existing row 159 ("Feria Lagunitas de Puerto Montt" / Mandarina data),
# pushing all subsequent rows (old 159..183) down by one (to 160..184),
# and the new row 159 is populated with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 159; everything from old row 159 downward
# shifts down to make room (old 159 -> 160, ..., old 183 -> 184).
$ws.Rows.Item(159).Insert()

# Populate the new row 159 with the new weekly observation.
$ws.Range("A159").Value = 4
$ws.Range("B159").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C159").Value = "Los Lagos"
$ws.Range("D159").Value = 44637
$ws.Range("E159").Value = 10
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100102
$ws.Range("H159").Value = "Cítricos"
$ws.Range("I159").Value = 100102004
$ws.Range("J159").Value = "Mandarina"
$ws.Range("K159").Value = "Murcott"
$ws.Range("L159").Value = "Primera"
$ws.Range("M159").Value = 500
$ws.Range("N159").Value = 12500
$ws.Range("O159").Value = 13000
$ws.Range("P159").Value = 12750
$ws.Range("Q159").Value = "$/bandeja 10 kilos"
$ws.Range("R159").Value = "Región de O'Higgins"
$ws.Range("S159").Value = 1275
$ws.Range("T159").Value = 10
